$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.912.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.875.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7417"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'242.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.9987"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3158"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.77%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07188"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'24.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.09%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08391"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.55%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.7508"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.98%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.419"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.25%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.875.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.50%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'92.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.94%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'29.901.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'6.096"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'243.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.93%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007813"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.77%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.9982"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.28%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.122.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.989"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.30%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.9993"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.16%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1554"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.280"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'165.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.35%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.511"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.30%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.592"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.69%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.531"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.274"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.59%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05317"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.79%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7541"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9969"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.71%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.753"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.4529"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.46%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.112.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.051"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8584"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.65%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'103.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.19%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Aptos"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'7.650"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'SynthetixNetwork"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'3.097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.31%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.840"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.020.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -8.19%  "
$ws.Range("E51").Style = "Normal"
